$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last 4 rows (rows 10-13), which correspond to the dropped "ECs" target
# cluster for every source cluster, collapsing the table from 12 to 8 data rows.
$ws.Range("A10:A13").EntireRow.Delete()

# Re-populate rows 2-9 with the recalculated TPM-derived values and the new
# (source cluster -> target cluster) pairing (target cluster "ECs" removed).

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ceacam1"
$ws.Range("C2").Value = "Havcr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 53.37999233333334
$ws.Range("H2").Value = 160.139977
$ws.Range("I2").Value = 0.9228671193042202
$ws.Range("J2").Value = 0.9228671193042202
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.088104
$ws.Range("N2").Value = 0.264312
$ws.Range("O2").Value = 0.00497951145580401
$ws.Range("P2").Value = 0.004979511455804011
$ws.Range("Q2").Value = 4.702990844536001
$ws.Range("R2").Value = 42.326917600824
$ws.Range("S2").Value = 0.00459542739276021
$ws.Range("T2").Value = 0.004595427392760211

# Row 3: ECs -> Resolving-Mac
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ceacam1"
$ws.Range("C3").Value = "Havcr2"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 53.37999233333334
$ws.Range("H3").Value = 160.139977
$ws.Range("I3").Value = 0.9228671193042202
$ws.Range("J3").Value = 0.9228671193042202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.605198
$ws.Range("N3").Value = 52.815594
$ws.Range("O3").Value = 0.9950204885441959
$ws.Range("P3").Value = 0.995020488544196
$ws.Range("Q3").Value = 939.7653342668153
$ws.Range("R3").Value = 8457.888008401338
$ws.Range("S3").Value = 0.9182716919114599
$ws.Range("T3").Value = 0.91827169191146

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ceacam1"
$ws.Range("C4").Value = "Havcr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1042846666666667
$ws.Range("H4").Value = 0.312854
$ws.Range("I4").Value = 0.001802939373113576
$ws.Range("J4").Value = 0.001802939373113576
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.088104
$ws.Range("N4").Value = 0.264312
$ws.Range("O4").Value = 0.00497951145580401
$ws.Range("P4").Value = 0.004979511455804011
$ws.Range("Q4").Value = 0.009187896272000002
$ws.Range("R4").Value = 0.08269106644800001
$ws.Range("S4").Value = 0.000008977757262539151
$ws.Range("T4").Value = 0.000008977757262539153

# Row 5: FAPs -> Resolving-Mac
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ceacam1"
$ws.Range("C5").Value = "Havcr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1042846666666667
$ws.Range("H5").Value = 0.312854
$ws.Range("I5").Value = 0.001802939373113576
$ws.Range("J5").Value = 0.001802939373113576
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.605198
$ws.Range("N5").Value = 52.815594
$ws.Range("O5").Value = 0.9950204885441959
$ws.Range("P5").Value = 0.995020488544196
$ws.Range("Q5").Value = 1.835952205030667
$ws.Range("R5").Value = 16.523569845276
$ws.Range("S5").Value = 0.001793961615851037
$ws.Range("T5").Value = 0.001793961615851037

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ceacam1"
$ws.Range("C6").Value = "Havcr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4198543333333333
$ws.Range("H6").Value = 1.259563
$ws.Range("I6").Value = 0.007258707657939662
$ws.Range("J6").Value = 0.007258707657939662
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.088104
$ws.Range("N6").Value = 0.264312
$ws.Range("O6").Value = 0.00497951145580401
$ws.Range("P6").Value = 0.004979511455804011
$ws.Range("Q6").Value = 0.036990846184
$ws.Range("R6").Value = 0.332917615656
$ws.Range("S6").Value = 0.00003614481793704284
$ws.Range("T6").Value = 0.00003614481793704284

# Row 7: MuSCs -> Resolving-Mac
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ceacam1"
$ws.Range("C7").Value = "Havcr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4198543333333333
$ws.Range("H7").Value = 1.259563
$ws.Range("I7").Value = 0.007258707657939662
$ws.Range("J7").Value = 0.007258707657939662
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.605198
$ws.Range("N7").Value = 52.815594
$ws.Range("O7").Value = 0.9950204885441959
$ws.Range("P7").Value = 0.995020488544196
$ws.Range("Q7").Value = 7.391618669491332
$ws.Range("R7").Value = 66.524568025422
$ws.Range("S7").Value = 0.007222562840002618
$ws.Range("T7").Value = 0.007222562840002619

# Row 8: Resolving-Mac -> FAPs
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ceacam1"
$ws.Range("C8").Value = "Havcr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.937340333333333
$ws.Range("H8").Value = 11.812021
$ws.Range("I8").Value = 0.06807123366472666
$ws.Range("J8").Value = 0.06807123366472666
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.088104
$ws.Range("N8").Value = 0.264312
$ws.Range("O8").Value = 0.00497951145580401
$ws.Range("P8").Value = 0.004979511455804011
$ws.Range("Q8").Value = 0.346895432728
$ws.Range("R8").Value = 3.122058894552
$ws.Range("S8").Value = 0.000338961487844218
$ws.Range("T8").Value = 0.000338961487844218

# Row 9: Resolving-Mac -> Resolving-Mac
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ceacam1"
$ws.Range("C9").Value = "Havcr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.937340333333333
$ws.Range("H9").Value = 11.812021
$ws.Range("I9").Value = 0.06807123366472666
$ws.Range("J9").Value = 0.06807123366472666
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.605198
$ws.Range("N9").Value = 52.815594
$ws.Range("O9").Value = 0.9950204885441959
$ws.Range("P9").Value = 0.995020488544196
$ws.Range("Q9").Value = 69.31765616171933
$ws.Range("R9").Value = 623.858905455474
$ws.Range("S9").Value = 0.06773227217688244
$ws.Range("T9").Value = 0.06773227217688244
